# Updates cryptos list values per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.896.51"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "2.625.80"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'520.82"
$ws.Range("E5").Value = "  +2.30%  "
$ws.Range("D6").Value = "'145.19"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "2.638.33"
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "3.089.86"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "58.914.31"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "'20.86"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "2.635.86"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("D19").Value = "'344.89"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'4.48"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "'10.21"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'61.69"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").Value = "'0.166"
$ws.Range("E26").Value = "  +3.50%  "
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("D29").Value = "'7.11"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "'6.24"
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'18.88"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.57"
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("D34").Value = "'150.45"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").Value = "'0.979"
$ws.Range("E35").Value = "  +3.34%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").Value = "'36.65"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("D39").Value = "'0.840"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("D40").Value = "'3.65"
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("E41").Value = "  +1.42%  "
$ws.Range("D42").Value = "'277.59"
$ws.Range("E42").Value = "  -5.22%  "
$ws.Range("D43").Value = "'0.996"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "'0.0985"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D46").Value = "'19.49"
$ws.Range("E46").Value = "  +2.18%  "
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("D48").Value = "'10.28"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "1.988.82"
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("D50").Value = "'0.0228"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Value = "'4.64"
$ws.Range("E51").Value = "  -1.14%  "
